$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.771.59"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.279.37"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "120.79"
$ws.Range("E5").Value = "  +6.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.99"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +5.00%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.50"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.34"
$ws.Range("E12").Value = "  +6.49%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.70"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.917"
$ws.Range("E15").Value = "  +6.95%  "
$ws.Range("D16").Value = "2.622.60"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "2.271.84"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "43.721.36"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.34"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.02"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.52"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.88"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.96"
$ws.Range("E26").Value = "  +5.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "43.38"
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.40"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.24"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.72"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0927"
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.79"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("E36").Value = "  +12.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0387"
$ws.Range("E37").Value = "  +10.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("E40").Value = "  +5.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.15"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.92"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.40"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "75.01"
$ws.Range("E47").Value = "  +43.71%  "
$ws.Range("E48").Value = "  +4.03%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.13"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.50"
$ws.Range("E51").Value = "  -1.24%  "
